$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.601.75'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.442.87'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.86'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.78'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').Value = '2.439.25'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  -4.04%  '
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.07'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('E15').Value = '  -5.23%  '
$ws.Range('D16').Value = '2.887.00'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '62.404.93'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').Value = '2.436.19'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.21'
$ws.Range('E19').Value = '  -3.16%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '327.04'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.17'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.10'
$ws.Range('E23').Value = '  +11.04%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.22'
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '626.54'
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.07'
$ws.Range('E27').Value = '  +3.30%  '
$ws.Range('D28').Value = '0.0₃0994'
$ws.Range('E28').Value = '  -5.92%  '
$ws.Range('D29').Value = '2.557.82'
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.15'
$ws.Range('E32').Value = '  -4.68%  '
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.138'
$ws.Range('E34').Value = '  -4.33%  '
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  -3.17%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.80'
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '146.63'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  -4.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.58'
$ws.Range('E43').Value = '  -3.33%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '146.78'
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.70'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('E48').Value = '  -4.07%  '
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0920'
$ws.Range('E51').Value = '  -1.14%  '
